$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.053496360778809
$ws.Range("B1").Value = 2.255256175994873
$ws.Range("C1").Value = 2.268766164779663
$ws.Range("D1").Value = 2.824041128158569
$ws.Range("E1").Value = 3.568579196929932
